# Applies the "submission obj & reb" edit:
#  - Sets the Active (column J) flag to FALSE for rows 2-13 on the Regression sheet
#  - Updates the active cell selection to J3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Regression")

# Set J2:J13 (Active column) values to FALSE
$ws.Range("J2:J13").Value = $false

# Activate the sheet and move the selection to J3
$ws.Activate()
$ws.Range("J3").Select()
